$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.951.49"
$ws.Range("D3").Value = "2.240.49"
$ws.Range("E3").Value = "  +2.04%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.10"
$ws.Range("E5").Value = "  +2.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.624"
$ws.Range("E6").Value = "  -0.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.70"
$ws.Range("E7").Value = "  -2.60%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +2.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.97"
$ws.Range("E10").Value = "  +1.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0921"
$ws.Range("E11").Value = "  +7.20%  "
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("D13").Value = "2.571.60"
$ws.Range("E13").Value = "  +2.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.78"
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.50"
$ws.Range("E15").Value = "  +2.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.809"
$ws.Range("E16").Value = "  -0.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.63"
$ws.Range("E17").Value = "  +1.45%  "
$ws.Range("D18").Value = "2.250.25"
$ws.Range("E18").Value = "  +3.13%  "
$ws.Range("D19").Value = "41.893.29"
$ws.Range("E19").Value = "  +4.81%  "
$ws.Range("D20").Value = "0.0₃0920"
$ws.Range("E20").Value = "  +1.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.71"
$ws.Range("E21").Value = "  +0.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.09"
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.01"
$ws.Range("E23").Value = "  +8.72%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.41"
$ws.Range("E25").Value = "  +2.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.31"
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("E27").Value = "  +1.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.146"
$ws.Range("E28").Value = "  +3.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "169.74"
$ws.Range("E29").Value = "  -1.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.15"
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("E31").Value = "  -1.38%  "
$ws.Range("E32").Value = "  -0.53%  "
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("E34").Value = "  +7.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.72"
$ws.Range("E35").Value = "  +3.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0641"
$ws.Range("E36").Value = "  +2.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.69"
$ws.Range("E37").Value = "  -4.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.74"
$ws.Range("E38").Value = "  -3.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.38"
$ws.Range("E39").Value = "  -2.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.000258"
$ws.Range("E40").Value = "  +29.76%  "
$ws.Range("E41").Value = "  +0.21%  "
$ws.Range("E42").Value = "  +5.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.69"
$ws.Range("E43").Value = "  +4.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.62"
$ws.Range("E44").Value = "  -8.44%  "
$ws.Range("E45").Value = "  +0.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.00"
$ws.Range("E46").Value = "  -2.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0960"
$ws.Range("E47").Value = "  +3.34%  "
$ws.Range("D48").Value = "1.488.16"
$ws.Range("E48").Value = "  -1.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.75"
$ws.Range("E49").Value = "  -3.99%  "
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("E51").Value = "  +5.09%  "
